$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy"
$ws.Range("G3").Value = "Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi"
$ws.Range("G4").Value = "Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid"
$ws.Range("G5").Value = "Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad"
$ws.Range("G6").Value = "Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub"
$ws.Range("G7").Value = "Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Rana Abo-Zaid"
$ws.Range("G8").Value = "Administrator, Dr. Asmaa Reda, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Majorelle Magdy"
$ws.Range("G9").Value = "Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Majorelle Magdy"
$ws.Range("G10").Value = "Dr. Heba Mahmoud Ali, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad"
$ws.Range("G11").Value = "Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat"
$ws.Range("G12").Value = "Dr. Salma El-Gendy, Administrator"
$ws.Range("G13").Value = "Dr. Safa Hany, Dr. Shimaa Ashraf, D Wessam Atef, Dr. Mariam Nour El-Din, Dr. Omnia Mohammad"
$ws.Range("G14").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G16").Value = "Dr. Nourhan Mohammad, Dr. Amal Awwad"
$ws.Range("G17").Value = "Dr. Nourhan Osama, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Dina Adel"
$ws.Range("G19").Value = "Dr. Sarah Mahdy, D Mariam E. Mohammad"
$ws.Range("G23").Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Range("G24").Value = "Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Maryam Ashraf, Dr. Salma Hassan, Dr. Monica, Dr. Remon, Dr. Youstina Magdy, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Aya Emad, Dr. Wafaa Ebida"
$ws.Range("G25").Value = "Dr. Marina Atef, Dr. Remon, Dr. Youstina Magdy, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Aya Emad"
$ws.Range("G26").Value = "Dr. Gehad Salah, Dr. Youstina Magdy"
$ws.Range("G27").Value = "Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Eman Mohammad Al"
$ws.Range("G28").Value = "Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida"
$ws.Range("G29").Value = "Dr. Neveen Nashaat, Dr. Monica, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Naema Gomaa"
$ws.Range("G30").Value = "Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid"
$ws.Range("G31").Value = "Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi"
$ws.Range("G32").Value = "Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid"
$ws.Range("G33").Value = "Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad"
$ws.Range("G34").Value = "Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub"
$ws.Range("G35").Value = "Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Rana Abo-Zaid"
$ws.Range("G36").Value = "Administrator, Dr. Asmaa Reda, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Majorelle Magdy"
$ws.Range("G37").Value = "Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Majorelle Magdy"
$ws.Range("G38").Value = "Dr. Heba Mahmoud Ali, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad"
$ws.Range("G39").Value = "Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat"
$ws.Range("G40").Value = "Dr. Salma El-Gendy, Administrator"
$ws.Range("G41").Value = "Dr. Safa Hany, Dr. Shimaa Ashraf, D Wessam Atef, Dr. Mariam Nour El-Din, Dr. Omnia Mohammad"
$ws.Range("G42").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G44").Value = "Dr. Nourhan Mohammad, Dr. Amal Awwad"
$ws.Range("G45").Value = "Dr. Nourhan Osama, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Dina Adel"
$ws.Range("G47").Value = "Dr. Sarah Mahdy, D Mariam E. Mohammad"
$ws.Range("G51").Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Range("G52").Value = "Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Maryam Ashraf, Dr. Salma Hassan, Dr. Monica, Dr. Remon, Dr. Youstina Magdy, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Aya Emad, Dr. Wafaa Ebida"
$ws.Range("G53").Value = "Dr. Marina Atef, Dr. Remon, Dr. Youstina Magdy, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Aya Emad"
$ws.Range("G54").Value = "Dr. Gehad Salah, Dr. Youstina Magdy"
$ws.Range("G55").Value = "Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Eman Mohammad Al"
$ws.Range("G56").Value = "Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida"
$ws.Range("G57").Value = "Dr. Neveen Nashaat, Dr. Monica, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Naema Gomaa"
